$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Row 35: 79. Word Search ----
$ws.Range("A35").Value = "79. Word Search"
$ws.Range("B35").Value = "Medium"
$ws.Range("C35").Value = "Backtracking"
$ws.Range("D35").Value = "Recursive Backtracking - DFS. Use a Set to not revisit the same position twice in the path. Use 2d grid traversal for pathfinding, call dfs on each 4 directions, but remove the cell after (backtracking)."
$ws.Range("E35").Value = "https://leetcode.com/problems/word-search/solutions/279279/java-simple-with-explanation/ "
$ws.Hyperlinks.Add($ws.Range("E35"), "https://leetcode.com/problems/word-search/solutions/279279/java-simple-with-explanation/ ") | Out-Null

# ---- Row 36: 322. Coin Change ----
$ws.Range("A36").Value = "322. Coin Change"
$ws.Range("B36").Value = "Medium"
$ws.Range("C36").Value = "Dynamic Programming"
$ws.Range("D36").Value = "dp[0] to dp[amount]. Track the min amount for each value in the range. Arrays.fill() the initial values to amount + 1, and only return at the end if the value does not equal the initial value. Start with the DFS - Backtracking approach and evolve to Top-Down Memoization, then DP Bottom-Up."
$ws.Range("E36").Value = "https://leetcode.com/problems/coin-change/solutions/778548/c-dp-solution-explained-100-time-100-space/ "
$ws.Hyperlinks.Add($ws.Range("E36"), "https://leetcode.com/problems/coin-change/solutions/778548/c-dp-solution-explained-100-time-100-space/ ") | Out-Null

# ---- Match formatting of the existing table rows (fill on Difficulty, hyperlink style on Link) ----
$ws.Range("B35").Interior.Color = $ws.Range("B33").Interior.Color
$ws.Range("B36").Interior.Color = $ws.Range("B33").Interior.Color

$ws.Range("E34").Copy() | Out-Null
$ws.Range("E35").PasteSpecial(-4122) | Out-Null
$ws.Range("E34").Copy() | Out-Null
$ws.Range("E36").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---- Expand the table / autofilter to include the new rows ----
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E36"))

# ---- Update view state to match the saved selection/scroll position ----
$ws.Range("D40").Select()
$excel.ActiveWindow.ScrollRow = 13
